$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. Copy number formats from the (now-shifted) reference column F into the
#    two new columns D and E, row-range by row-range (skipping label-only
#    rows 36/37 and 78/79 which must stay untouched).
$fmtRanges = @("7:35", "38:77", "80:102")
foreach ($rr in $fmtRanges) {
    $bounds = $rr.Split(":")
    $r1 = $bounds[0]
    $r2 = $bounds[1]
    $ws.Range("F" + $r1 + ":F" + $r2).Copy()
    $ws.Range("D" + $r1 + ":D" + $r2).PasteSpecial(-4122)
    $ws.Range("E" + $r1 + ":E" + $r2).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# 3. Populate the two new columns with the new quarter's figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 2224000
$ws.Range("E8").Value2 = 2298900
$ws.Range("D9").Value2 = 1544000
$ws.Range("E9").Value2 = 1565800
$ws.Range("D10").Value2 = 680000
$ws.Range("E10").Value2 = 733100
$ws.Range("D12").Value2 = 15000
$ws.Range("E12").Value2 = 20400
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 29000
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 1769000
$ws.Range("E17").Value2 = 1765200
$ws.Range("D18").Value2 = 455000
$ws.Range("E18").Value2 = 533700
$ws.Range("D20").Value2 = 74000
$ws.Range("E20").Value2 = 29400
$ws.Range("D21").Value2 = 787000
$ws.Range("E21").Value2 = 820300
$ws.Range("D22").Value2 = 39900
$ws.Range("E22").Value2 = 34200
$ws.Range("D23").Value2 = 489100
$ws.Range("E23").Value2 = 528900
$ws.Range("D24").Value2 = 147700
$ws.Range("E24").Value2 = 38900
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 341400
$ws.Range("E26").Value2 = 490000
$ws.Range("D27").Value2 = 331900
$ws.Range("E27").Value2 = 483200
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 15600
$ws.Range("E29").Value2 = -30300
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -74000
$ws.Range("E32").Value2 = -29400
$ws.Range("D33").Value2 = 347500
$ws.Range("E33").Value2 = 452900
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 347500
$ws.Range("E35").Value2 = 452900
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 2923300
$ws.Range("E41").Value2 = 2791300
$ws.Range("D42").Value2 = 12300
$ws.Range("E42").Value2 = 184700
$ws.Range("D43").Value2 = 1613800
$ws.Range("E43").Value2 = 1441900
$ws.Range("D44").Value2 = 403400
$ws.Range("E44").Value2 = 473600
$ws.Range("D45").Value2 = 137100
$ws.Range("E45").Value2 = 190700
$ws.Range("D46").Value2 = 5089900
$ws.Range("E46").Value2 = 5082200
$ws.Range("D47").Value2 = 2228300
$ws.Range("E47").Value2 = 2290500
$ws.Range("D48").Value2 = 9959800
$ws.Range("E48").Value2 = 9923700
$ws.Range("D49").Value2 = 1197300
$ws.Range("E49").Value2 = 1227400
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 666700
$ws.Range("E52").Value2 = 654500
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 19142000
$ws.Range("E54").Value2 = 19178300
$ws.Range("D57").Value2 = 1708000
$ws.Range("E57").Value2 = 1763600
$ws.Range("D58").Value2 = 453300
$ws.Range("E58").Value2 = 460900
$ws.Range("D59").Value2 = 142200
$ws.Range("E59").Value2 = 113800
$ws.Range("D60").Value2 = 2303500
$ws.Range("E60").Value2 = 2338300
$ws.Range("D61").Value2 = 3314600
$ws.Range("E61").Value2 = 3351700
$ws.Range("D62").Value2 = 2320500
$ws.Range("E62").Value2 = 2312000
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 8259100
$ws.Range("E66").Value2 = 8320800
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 13497900
$ws.Range("E72").Value2 = 13409900
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 10882900
$ws.Range("E76").Value2 = 10857500
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 347500
$ws.Range("E81").Value2 = 452900
$ws.Range("D83").Value2 = 258000
$ws.Range("E83").Value2 = 257200
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 655200
$ws.Range("E89").Value2 = 691400
$ws.Range("D91").Value2 = -403400
$ws.Range("E91").Value2 = -410300
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -226500
$ws.Range("E94").Value2 = -610200
$ws.Range("D96").Value2 = -241500
$ws.Range("E96").Value2 = -241200
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -289800
$ws.Range("E100").Value2 = -262000
$ws.Range("D101").Value2 = -6900
$ws.Range("E101").Value2 = -14400
$ws.Range("D102").Value2 = 132000
$ws.Range("E102").Value2 = -195200

